# Update leve profit calculations across multiple crafting-job sheets
# (ALC, ARM, BSM, CRP, CUL, WVR) - refreshed market-board price data
# pulled in by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 18519048
$ws.Range("I19").Value = 608.8889
$ws.Range("J19").Value = 27778268
$ws.Range("K19").Value = 608.8889
$ws.Range("L19").Value = 27778268
$ws.Range("M19").Value = -433.8889
$ws.Range("N19").Value = -27778618
$ws.Range("H63").Value = 47475
$ws.Range("J63").Value = 47475
$ws.Range("L63").Value = 47475
$ws.Range("N63").Value = -48723
$ws.Range("H66").Value = 47475
$ws.Range("J66").Value = 47475
$ws.Range("L66").Value = 142425
$ws.Range("N66").Value = -148665
$ws.Range("H70").Value = 1250.069
$ws.Range("I70").Value = 1266.2307
$ws.Range("J70").Value = 1110
$ws.Range("K70").Value = 3798.6921
$ws.Range("L70").Value = 3330
$ws.Range("M70").Value = -3528.6921
$ws.Range("N70").Value = -3870
$ws.Range("H73").Value = 1250.069
$ws.Range("I73").Value = 1266.2307
$ws.Range("J73").Value = 1110
$ws.Range("K73").Value = 3798.6921
$ws.Range("L73").Value = 3330
$ws.Range("M73").Value = -2862.6921
$ws.Range("N73").Value = -5202
$ws.Range("H74").Value = 4210.4443
$ws.Range("I74").Value = 3679
$ws.Range("J74").Value = 5273.3335
$ws.Range("K74").Value = 3679
$ws.Range("L74").Value = 5273.3335
$ws.Range("M74").Value = -2743
$ws.Range("N74").Value = -7145.3335
$ws.Range("H77").Value = 4210.4443
$ws.Range("I77").Value = 3679
$ws.Range("J77").Value = 5273.3335
$ws.Range("K77").Value = 18395
$ws.Range("L77").Value = 26366.6675
$ws.Range("M77").Value = -13715
$ws.Range("N77").Value = -35726.6675
$ws.Range("H132").Value = 4940.095
$ws.Range("I132").Value = 3962.3403
$ws.Range("J132").Value = 7812.25
$ws.Range("K132").Value = 11887.0209
$ws.Range("L132").Value = 23436.75
$ws.Range("M132").Value = -9357.0209
$ws.Range("N132").Value = -28496.75
$ws.Range("H138").Value = 2194.4126
$ws.Range("I138").Value = 1425.8636
$ws.Range("J138").Value = 3133.75
$ws.Range("K138").Value = 4277.5908
$ws.Range("L138").Value = 9401.25
$ws.Range("M138").Value = 862.4092000000001
$ws.Range("N138").Value = -19681.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 800
$ws.Range("I2").Value = 600
$ws.Range("K2").Value = 600
$ws.Range("M2").Value = -487
$ws.Range("H32").Value = 15820.947
$ws.Range("I32").Value = 8922.338
$ws.Range("J32").Value = 45331.668
$ws.Range("K32").Value = 8922.338
$ws.Range("L32").Value = 45331.668
$ws.Range("M32").Value = -8635.338
$ws.Range("N32").Value = -45905.668
$ws.Range("H45").Value = 887.1579
$ws.Range("I45").Value = 889.4666999999999
$ws.Range("J45").Value = 878.5
$ws.Range("K45").Value = 889.4666999999999
$ws.Range("L45").Value = 878.5
$ws.Range("M45").Value = -512.4666999999999
$ws.Range("N45").Value = -1632.5
$ws.Range("H54").Value = 20049
$ws.Range("J54").Value = 20049
$ws.Range("L54").Value = 20049
$ws.Range("N54").Value = -21587
$ws.Range("H62").Value = 19250
$ws.Range("J62").Value = 19250
$ws.Range("L62").Value = 19250
$ws.Range("N62").Value = -20498
$ws.Range("H65").Value = 19250
$ws.Range("J65").Value = 19250
$ws.Range("L65").Value = 57750
$ws.Range("N65").Value = -63990
$ws.Range("H103").Value = 26180.5
$ws.Range("J103").Value = 26180.5
$ws.Range("L103").Value = 26180.5
$ws.Range("N103").Value = -28524.5
$ws.Range("H106").Value = 39800
$ws.Range("J106").Value = 39800
$ws.Range("L106").Value = 39800
$ws.Range("N106").Value = -42324
$ws.Range("H116").Value = 800
$ws.Range("I116").Value = 600
$ws.Range("K116").Value = 600
$ws.Range("M116").Value = 1694

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 800
$ws.Range("I3").Value = 600
$ws.Range("K3").Value = 600
$ws.Range("M3").Value = -486

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1923.6522
$ws.Range("I134").Value = 1372
$ws.Range("K134").Value = 4116
$ws.Range("M134").Value = -1581

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H43").Value = 5880.2
$ws.Range("I43").Value = 3002
$ws.Range("J43").Value = 6200
$ws.Range("K43").Value = 9006
$ws.Range("L43").Value = 18600
$ws.Range("M43").Value = -8892
$ws.Range("N43").Value = -18828
$ws.Range("H44").Value = 327
$ws.Range("I44").Value = 311.44446
$ws.Range("J44").Value = 397
$ws.Range("K44").Value = 934.33338
$ws.Range("L44").Value = 1191
$ws.Range("M44").Value = -536.33338
$ws.Range("N44").Value = -1987
$ws.Range("H64").Value = 8508.666999999999
$ws.Range("I64").Value = 512
$ws.Range("J64").Value = 12507
$ws.Range("K64").Value = 1536
$ws.Range("L64").Value = 37521
$ws.Range("M64").Value = -1266
$ws.Range("N64").Value = -38061
$ws.Range("H67").Value = 8508.666999999999
$ws.Range("I67").Value = 512
$ws.Range("J67").Value = 12507
$ws.Range("K67").Value = 1536
$ws.Range("L67").Value = 37521
$ws.Range("M67").Value = -600
$ws.Range("N67").Value = -39393
$ws.Range("H114").Value = 1460.2667
$ws.Range("I114").Value = 414.85715
$ws.Range("J114").Value = 2375
$ws.Range("K114").Value = 1244.57145
$ws.Range("L114").Value = 7125
$ws.Range("M114").Value = 2009.42855
$ws.Range("N114").Value = -13633
$ws.Range("H117").Value = 20411930
$ws.Range("I117").Value = 1286
$ws.Range("J117").Value = 35719916
$ws.Range("K117").Value = 3858
$ws.Range("L117").Value = 107159748
$ws.Range("M117").Value = -416
$ws.Range("N117").Value = -107166632
$ws.Range("H134").Value = 2573.7437
$ws.Range("J134").Value = 3826.9524
$ws.Range("L134").Value = 11480.8572
$ws.Range("N134").Value = -21620.8572

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 8000
$ws.Range("J47").Value = 8000
$ws.Range("L47").Value = 8000
$ws.Range("N47").Value = -9144
$ws.Range("H69").Value = 8127
$ws.Range("J69").Value = 8127
$ws.Range("L69").Value = 8127
$ws.Range("N69").Value = -9625
$ws.Range("H72").Value = 8127
$ws.Range("J72").Value = 8127
$ws.Range("L72").Value = 24381
$ws.Range("N72").Value = -31869
$ws.Range("H101").Value = 13551
$ws.Range("J101").Value = 13551
$ws.Range("L101").Value = 13551
$ws.Range("N101").Value = -20041
$ws.Range("H104").Value = 8545
$ws.Range("J104").Value = 8545
$ws.Range("L104").Value = 8545
$ws.Range("N104").Value = -15533
